# Applies the commit "Added readme and script to process antibiotic
# resistance data" to the workbook.
#
# Net effect observed in the target OOXML:
#  - Sheet2 worksheet (physically xl/worksheets/sheet1.xml): the B/CARB
#    row (row 32) and part of the D/KAN row (row 34, cols H:J) were
#    overwritten with values that belong to the "Carb"/"Kan" mini summary
#    table living at C41:L44 on the same sheet (cols D:I of that table).
#    The summary table's own cached values (rows 42-43) were also
#    refreshed (same numbers, refreshed float precision).
#  - Sheet1 worksheet (physically xl/worksheets/sheet2.xml): the
#    CARB row (row 32, cols E:G) got new values, and an entirely new
#    "Carb"/"Kan" mini summary table (mirroring the one on Sheet2) was
#    added at E45:N48.
#  - Sheet3 worksheet gained no cell-value changes; only selection/active
#    tab bookkeeping changed (Sheet3 becomes the active tab/sheet).

$wb  = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws3 = $wb.Worksheets.Item("Sheet3")

# ---------------------------------------------------------------------
# Sheet2: CARB row (32) + KAN row (34, H:J) refreshed from the summary
# table, and the summary table itself (rows 42-43) refreshed.
# ---------------------------------------------------------------------
$ws2.Range("B32").Value = 1.6219999790191699
$ws2.Range("C32").Value = 0.36669999361038202
$ws2.Range("D32").Value = 0.3817999958992
$ws2.Range("E32").Value = 0.047100000083446503
$ws2.Range("F32").Value = 0.046000000089407002
$ws2.Range("G32").Value = 0.0472000017762184

$ws2.Range("H34").Value = 0.0465999990701675
$ws2.Range("I34").Value = 0.046700000762939502
$ws2.Range("J34").Value = 0.339300006628036

$ws2.Range("D42").Value = 1.6219999790191699
$ws2.Range("E42").Value = 0.36669999361038202
$ws2.Range("F42").Value = 0.3817999958992
$ws2.Range("H42").Value = 0.046000000089407002
$ws2.Range("I42").Value = 0.0472000017762184

$ws2.Range("D43").Value = 0.0465999990701675
$ws2.Range("E43").Value = 0.046700000762939502
$ws2.Range("F43").Value = 0.339300006628036

# ---------------------------------------------------------------------
# Sheet1: CARB row (32, E:G) refreshed.
# ---------------------------------------------------------------------
$ws1.Range("E32").Value = 0.048799999058246599
$ws1.Range("F32").Value = 1.2474999427795399
$ws1.Range("G32").Value = 1.50960004329681

# ---------------------------------------------------------------------
# Sheet1: new Carb/Kan mini summary table at E45:N48, mirroring the
# existing one on Sheet2 (C41:L44). Copy it over first (to pick up
# labels/styling/shared strings) then overwrite the refreshed numbers.
# ---------------------------------------------------------------------
$ws2.Range("B42").Copy($ws1.Range("D46"))
$ws2.Range("C41:L41").Copy($ws1.Range("E45"))
$ws2.Range("C42:L42").Copy($ws1.Range("E46"))
$ws2.Range("C43:F43").Copy($ws1.Range("E47"))
$ws2.Range("D44:F44").Copy($ws1.Range("F48"))

$ws1.Range("F46").Value = 1.6219999790191699
$ws1.Range("G46").Value = 0.36669999361038202
$ws1.Range("H46").Value = 0.3817999958992
$ws1.Range("I46").Value = 0.047100000083446503
$ws1.Range("J46").Value = 0.046000000089407002
$ws1.Range("K46").Value = 0.0472000017762184
$ws1.Range("L46").Value = 0.048799999058246599
$ws1.Range("M46").Value = 1.2474999427795399
$ws1.Range("N46").Value = 1.50960004329681

$ws1.Range("F47").Value = 0.0465999990701675
$ws1.Range("G47").Value = 0.046700000762939502
$ws1.Range("H47").Value = 0.339300006628036

# ---------------------------------------------------------------------
# Selection / active-tab bookkeeping (Sheet3 ends up the active sheet,
# each sheet's selection moves to where the editor last left it).
# ---------------------------------------------------------------------
$ws2.Activate()
$ws2.Range("A13").Select()
$ws2.Range("H50").Select()

$ws1.Activate()
$ws1.Range("A22").Select()
$ws1.Range("J43").Select()

$ws3.Activate()
$ws3.Range("A21").Select()
$ws3.Range("I53").Select()
